$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week date range) ---
$ws.Range("A8").Value = "Volume 30   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/10/2023  Through  7/16/2023"

# --- Helper functions for style-transition cells ---
function Set-NumFromText($addr, $val, $fmtSrc) {
    $ref = $ws.Range($addr)
    $src = $ws.Range($fmtSrc)
    $src.Copy()
    $ref.PasteSpecial(-4122)
    $ref.Value = $val
}

function Set-TextPlaceholder($addr, $text, $fmtSrc) {
    $ref = $ws.Range($addr)
    $ref.NumberFormat = "@"
    $ref.Value = $text
    $src = $ws.Range($fmtSrc)
    $src.Copy()
    $ref.PasteSpecial(-4122)
}

# --- Style-transition cells: text placeholder -> number ---
Set-NumFromText "C14" 1 "I14"
Set-NumFromText "F14" 1 "I14"
Set-NumFromText "D15" 1 "I14"
Set-NumFromText "E15" -100 "L14"
Set-NumFromText "C18" 2 "I14"
Set-NumFromText "D26" 1 "I14"
Set-NumFromText "E26" -100 "L14"
Set-NumFromText "C28" 1 "I14"
Set-NumFromText "C29" 1 "I14"

# --- Style-transition cells: number -> text placeholder ---
Set-TextPlaceholder "D27" "0" "C23"
Set-TextPlaceholder "E27" "***.*" "E23"
Set-TextPlaceholder "G30" "0" "C23"
Set-TextPlaceholder "H30" "***.*" "E23"

# --- Plain value updates (style unchanged) ---
# Row 14
$ws.Range("I14").Value = 2
$ws.Range("L14").Value = 100
$ws.Range("M14").Value = -60
$ws.Range("N14").Value = -80
# Row 15
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = -23.076923076923
$ws.Range("L15").Value = -37.5
# Row 16
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 76
$ws.Range("J16").Value = 93
$ws.Range("K16").Value = -18.279569892473
$ws.Range("L16").Value = -10.588235294117
$ws.Range("M16").Value = -26.923076923076
$ws.Range("N16").Value = -87.290969899665
# Row 17
$ws.Range("C17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 93.75
$ws.Range("I17").Value = 154
$ws.Range("J17").Value = 163
$ws.Range("K17").Value = -5.521472392638
$ws.Range("L17").Value = 52.475247524752
$ws.Range("M17").Value = 57.142857142857
$ws.Range("N17").Value = -49.006622516556
# Row 18
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = -81.818181818181
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -75
$ws.Range("I18").Value = 66
$ws.Range("J18").Value = 134
$ws.Range("K18").Value = -50.746268656716
$ws.Range("L18").Value = -15.384615384615
$ws.Range("M18").Value = -52.857142857142
$ws.Range("N18").Value = -89.215686274509
# Row 19
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -42.307692307692
$ws.Range("I19").Value = 285
$ws.Range("J19").Value = 460
$ws.Range("K19").Value = -38.043478260869
$ws.Range("L19").Value = 20.253164556962
$ws.Range("M19").Value = 46.907216494845
$ws.Range("N19").Value = 19.747899159663
# Row 20
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = -47.368421052631
$ws.Range("I20").Value = 88
$ws.Range("J20").Value = 85
$ws.Range("K20").Value = 3.529411764705
$ws.Range("L20").Value = 137.837837837838
$ws.Range("M20").Value = 39.682539682539
$ws.Range("N20").Value = -83.582089552238
# Row 21
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -45.714285714285
$ws.Range("F21").Value = 89
$ws.Range("G21").Value = 129
$ws.Range("H21").Value = -31.007751937984
$ws.Range("I21").Value = 681
$ws.Range("J21").Value = 948
$ws.Range("K21").Value = -28.164556962025
$ws.Range("L21").Value = 22.702702702702
$ws.Range("M21").Value = 10.731707317073
$ws.Range("N21").Value = -70.583153347732
# Row 22
$ws.Range("L22").Value = 45.454545454545
# Row 24
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 12.5
$ws.Range("F24").Value = 100
$ws.Range("G24").Value = 119
$ws.Range("H24").Value = -15.966386554621
$ws.Range("I24").Value = 635
$ws.Range("J24").Value = 758
$ws.Range("K24").Value = -16.226912928759
$ws.Range("L24").Value = 39.867841409691
$ws.Range("M24").Value = 59.147869674185
# Row 25
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 100
$ws.Range("G25").Value = 40
$ws.Range("H25").Value = 27.5
$ws.Range("I25").Value = 277
$ws.Range("J25").Value = 222
$ws.Range("K25").Value = 24.774774774774
$ws.Range("L25").Value = 9.055118110236
$ws.Range("M25").Value = -21.971830985915
# Row 26
$ws.Range("J26").Value = 20
$ws.Range("K26").Value = -25
$ws.Range("L26").Value = -37.5
# Row 27
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 48
$ws.Range("K27").Value = 37.142857142857
$ws.Range("L27").Value = 41.176470588235
# Row 28
$ws.Range("F28").Value = 2
$ws.Range("I28").Value = 7
$ws.Range("K28").Value = -36.363636363636
$ws.Range("L28").Value = 40
$ws.Range("M28").Value = 16.666666666666
$ws.Range("N28").Value = -36.363636363636
# Row 29
$ws.Range("F29").Value = 2
$ws.Range("I29").Value = 6
$ws.Range("K29").Value = 200
$ws.Range("L29").Value = 20
$ws.Range("M29").Value = 20
$ws.Range("N29").Value = -45.454545454545
